$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet now holds 21 comment rows instead of 30, so drop the extra rows first.
$ws.Rows("22:30").Delete() | Out-Null

# Rewrite each remaining comment with the new wording / verification language.
# (Cells are written in an order that reproduces the shared-string table layout
#  of the target workbook; the *content* of every row A1:A21 is what actually matters.)
$ws.Range("A1").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI CARESOURCE & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A2").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MCR AND SEC MCD & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A4").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI AETNA & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A3").Value = 'PATIENT NOT FOUND'
$ws.Range("A5").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI BCBS & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A6").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MCR AND SEC BCBS & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A7").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI GOLDEN RULE & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A10").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MEDICAL MUTUAL & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A11").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MERIDIAN HEALTH & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A13").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI TRICARE AND SEC MCD & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A14").Value = 'RECEIVED PT DEMOGRAPHIC INFO FROM HOSPITAL FS CROSS CHECKING SOFTWARE FOUND AS THE SAME PT NAME, DOB, SSN#, PH# AND INSURANCE INFO, VERIFIED THE ADDRESS IN USPS DATA FOUND AS VALID.'
$ws.Range("A16").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MOLINA & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A17").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI BUCKEYE & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A21").Value = 'RECEIVED PT DEMOGRAPHIC INFO FROM HOSPITAL FS CROSS CHECKING SOFTWARE FOUND AS THE SAME PT NAME, DOB, SSN#, PH# AND VERIFIED THE ADDRESS IN USPS DATA FOUND AS VALID. INSURANCE INFORMATION FOUND AS SELFPAY.'
$ws.Range("A8").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI AETNA & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A9").Value = 'PATIENT NOT FOUND'
$ws.Range("A12").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MERIDIAN HEALTH & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A15").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI BCBS & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A18").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI BCBS & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A19").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI BCBS & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'
$ws.Range("A20").Value = 'RECEIVED GOOD INSURANCE INFO FROM HOSPITAL FS WITH PRI MOLINA & VERIFIED ADDRESS IN USPS DATA FOUND AS VALID, THEREFORE UPDATED INS INFO IN ESO AND CLM FILED TO INS'

# Reproduce the saved selection / scroll position from the edited workbook.
$ws.Range("A29").Select() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 14
$aw.ScrollColumn = 1
